# Generate Report for Handoff
# - Set Priority ("ht") for rows 7-12 on the zh-cn and de-de sheets
# - Refresh the "Latest Handoff Datetime" (column H) for those same rows

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

foreach ($row in 7..12) {
    $zhcn.Cells.Item($row, 5).Value = "ht"   # column E: Priority
    $dede.Cells.Item($row, 5).Value = "ht"   # column E: Priority
}

foreach ($row in 7..12) {
    $zhcn.Cells.Item($row, 8).Value = "2016-09-07 16:32:56"  # column H: Latest Handoff Datetime
    $dede.Cells.Item($row, 8).Value = "2016-09-07 16:33:09"  # column H: Latest Handoff Datetime
}
